$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Move the "X" value from L4 to M4 (cut preserves the shared-string reference)
$ws.Range("L4").Cut($ws.Range("M4"))
$ws.Range("L4").Clear()

# Update active selection to A4
$ws.Range("A4").Select()
